# Add a new paper entry to the literature review paper-bank table.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a new row at the end of the table; Word clones formatting from the
# row immediately above it (tcPr borders/margins, run/paragraph rPr, etc.)
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "Introduction to Airborne Disease Transmission Indoors"
$newRow.Cells.Item(2).Range.Text = "Importance and Impact"
$newRow.Cells.Item(3).Range.Text = "Poudel, 2021, " + [char]34 + "Impact of Covid-19 on health-related quality of life of patients: A structured review" + [char]34
